$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# ---------------------------------------------------------------------------
# 1) Text edits: split out "Quantity" into its own underlined run in the
#    MATERIAL(...) and PRODUCT(...) text boxes.
# ---------------------------------------------------------------------------

$material = $s.Shapes.Item(3)   # TextBox 4 -> "MATERIAL(MaterialCode, MaterialName, Quantity, Unit, MinLevel)"
$trMaterial = $material.TextFrame.TextRange
$qIdx = $trMaterial.Text.IndexOf("Quantity") + 1
$trMaterial.Characters($qIdx, 8).Font.Underline = 1

$product = $s.Shapes.Item(4)    # TextBox 6 -> "PRODUCT(ProductCode, ProductName, Description, Quantity, Unit, AvgCost)"
$trProduct = $product.TextFrame.TextRange
$qIdx2 = $trProduct.Text.IndexOf("Quantity") + 1
$trProduct.Characters($qIdx2, 8).Font.Underline = 1

# ---------------------------------------------------------------------------
# 2) Add six new connector shapes to the ER diagram, matching the style of
#    existing nearby connectors (dashed "line" connectors and triangle-arrow
#    "straightConnector1" connectors). We duplicate existing shapes so the
#    Quick-Style (p:style lnRef/fillRef/effectRef/fontRef) and line formatting
#    are carried over exactly, then move/resize/rename the duplicates.
# ---------------------------------------------------------------------------

$donorLine        = $s.Shapes.Item(15)  # "Straight Connector 32"       (line, lgDash, no flip)
$donorArrowFlipV   = $s.Shapes.Item(17) # "Straight Arrow Connector 39" (straightConnector1, flipV)
$donorArrowNoFlip  = $s.Shapes.Item(20) # "Straight Arrow Connector 45" (straightConnector1, no flip)

# The COM layer assigns each newly created shape the smallest shape-id not
# already used on the slide, and that cursor only ever advances (deleting a
# shape does not free its id for reuse). We burn through the unwanted ids
# with disposable placeholder connectors so the six kept duplicates land on
# the exact target ids (34, 36, 37, 39, 42, 47), then delete the throwaways.
$throwaways = New-Object System.Collections.ArrayList
function NewThrowaway() {
    $t = $s.Shapes.AddConnector(1, 1, 1, 2, 2)
    [void]$throwaways.Add($t.Name)
}

for ($k = 1; $k -le 17; $k++) { NewThrowaway }

$shape1 = $donorLine.Duplicate().Item(1)        # id 34 "Straight Connector 33"
$shape2 = $donorLine.Duplicate().Item(1)        # id 36 "Straight Connector 35"
$shape3 = $donorArrowFlipV.Duplicate().Item(1)  # id 37 "Straight Arrow Connector 36"

NewThrowaway

$shape4 = $donorArrowFlipV.Duplicate().Item(1)  # id 39 "Straight Arrow Connector 38"
$shape5 = $donorLine.Duplicate().Item(1)        # id 42 "Straight Connector 41"
$shape6 = $donorArrowNoFlip.Duplicate().Item(1) # id 47 "Straight Arrow Connector 46"

foreach ($nm in $throwaways) {
    $s.Shapes.Item($nm).Delete()
}

# Name + position/size each new shape (EMU / 12700 = points).
$shape1.Name = "Straight Connector 33"
$shape1.Left = 3444974 / 12700
$shape1.Top = 3155787 / 12700
$shape1.Width = 767262 / 12700
$shape1.Height = 0 / 12700

$shape2.Name = "Straight Connector 35"
$shape2.Left = 7185172 / 12700
$shape2.Top = 5751599 / 12700
$shape2.Width = 767262 / 12700
$shape2.Height = 0 / 12700

$shape3.Name = "Straight Arrow Connector 36"
$shape3.Left = 3759691 / 12700
$shape3.Top = 2703537 / 12700
$shape3.Width = 1484868 / 12700
$shape3.Height = 206191 / 12700

$shape4.Name = "Straight Arrow Connector 38"
$shape4.Left = 7608015 / 12700
$shape4.Top = 4653030 / 12700
$shape4.Width = 1710914 / 12700
$shape4.Height = 821204 / 12700

$shape5.Name = "Straight Connector 41"
$shape5.Left = 3215213 / 12700
$shape5.Top = 1993099 / 12700
$shape5.Width = 789611 / 12700
$shape5.Height = 9965 / 12700

$shape6.Name = "Straight Arrow Connector 46"
$shape6.Left = 3560783 / 12700
$shape6.Top = 2038345 / 12700
$shape6.Width = 1757964 / 12700
$shape6.Height = 340759 / 12700

# Move the six new shapes into the z-order right after "Straight Connector 32"
# (immediately before "Straight Connector 34"), preserving their relative order.
function MoveToPosition($shape, $targetPos) {
    while ($shape.ZOrderPosition -gt $targetPos) { $shape.ZOrder(3) }  # msoSendBackward
    while ($shape.ZOrderPosition -lt $targetPos) { $shape.ZOrder(4) }  # msoBringForward
}

$targetPos = $donorLine.ZOrderPosition + 1
MoveToPosition $shape1 $targetPos
$targetPos += 1
MoveToPosition $shape2 $targetPos
$targetPos += 1
MoveToPosition $shape3 $targetPos
$targetPos += 1
MoveToPosition $shape4 $targetPos
$targetPos += 1
MoveToPosition $shape5 $targetPos
$targetPos += 1
MoveToPosition $shape6 $targetPos
